$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value2 = 614.2857
$ws.Range("I9").Value2 = 166.66667
$ws.Range("J9").Value2 = 950
$ws.Range("K9").Value2 = 166.66667
$ws.Range("L9").Value2 = 950
$ws.Range("M9").Value2 = 2.333329999999989
$ws.Range("N9").Value2 = -1288

# Row 32
$ws.Range("H32").Value2 = 4132.3335
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 4132.3335
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 4132.3335
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value2 = -4784.3335

# Row 40
$ws.Range("H40").Value2 = 3958.75
$ws.Range("I40").Value2 = 3893.3333
$ws.Range("J40").Value2 = 3998
$ws.Range("K40").Value2 = 3893.3333
$ws.Range("L40").Value2 = 3998
$ws.Range("M40").Value2 = -3718.3333
$ws.Range("N40").Value2 = -4348

# Row 98
$ws.Range("H98").Value2 = 1430.4375
$ws.Range("I98").Value2 = 791.9286
$ws.Range("K98").Value2 = 791.9286
$ws.Range("M98").Value2 = 706.0714

# Row 122
$ws.Range("H122").Value2 = 1430.4375
$ws.Range("I122").Value2 = 791.9286
$ws.Range("K122").Value2 = 2375.7858
$ws.Range("M122").Value2 = 74.21420000000035

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 15118.31
$ws.Range("I32").Value2 = 14582.911
$ws.Range("K32").Value2 = 14582.911
$ws.Range("M32").Value2 = -14295.911

# Row 45
$ws.Range("H45").Value2 = 5834.0347
$ws.Range("I45").Value2 = 6677.7617
$ws.Range("J45").Value2 = 3619.25
$ws.Range("K45").Value2 = 6677.7617
$ws.Range("L45").Value2 = 3619.25
$ws.Range("M45").Value2 = -6300.7617
$ws.Range("N45").Value2 = -4373.25

# Row 74
$ws.Range("H74").Value2 = 37412.793
$ws.Range("I74").Value2 = 41306.152
$ws.Range("K74").Value2 = 41306.152
$ws.Range("M74").Value2 = -40432.152

# Row 77
$ws.Range("H77").Value2 = 37412.793
$ws.Range("I77").Value2 = 41306.152
$ws.Range("K77").Value2 = 206530.76
$ws.Range("M77").Value2 = -202162.76

# Row 102
$ws.Range("H102").Value2 = 2271.625
$ws.Range("I102").Value2 = 2235
$ws.Range("J102").Value2 = 2332.6667
$ws.Range("K102").Value2 = 2235
$ws.Range("L102").Value2 = 2332.6667
$ws.Range("M102").Value2 = -613
$ws.Range("N102").Value2 = -5576.6667

# Row 132
$ws.Range("H132").Value2 = 22450.604
$ws.Range("I132").Value2 = 24677.117
$ws.Range("J132").Value2 = 3302.6
$ws.Range("K132").Value2 = 74031.351
$ws.Range("L132").Value2 = 9907.799999999999
$ws.Range("M132").Value2 = -71501.351
$ws.Range("N132").Value2 = -14967.8

# Row 139
$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 3667.125
$ws.Range("J31").Value2 = 4989.3335
$ws.Range("L31").Value2 = 4989.3335
$ws.Range("N31").Value2 = -5579.3335

# Row 34
$ws.Range("H34").Value2 = 3667.125
$ws.Range("J34").Value2 = 4989.3335
$ws.Range("L34").Value2 = 4989.3335
$ws.Range("N34").Value2 = -5393.3335

# Row 94
$ws.Range("H94").Value2 = 2500
$ws.Range("I94").Value2 = 2500
$ws.Range("K94").Value2 = 2500
$ws.Range("M94").Value2 = -2049

# Row 122
$ws.Range("H122").Value2 = 1287.3125
$ws.Range("I122").Value2 = 1341.2858
$ws.Range("J122").Value2 = 909.5
$ws.Range("K122").Value2 = 4023.8574
$ws.Range("L122").Value2 = 2728.5
$ws.Range("M122").Value2 = -1573.8574
$ws.Range("N122").Value2 = -7628.5

# Row 132
$ws.Range("H132").Value2 = 2405.149
$ws.Range("I132").Value2 = 2200.2
$ws.Range("K132").Value2 = 6600.599999999999
$ws.Range("M132").Value2 = -4070.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 766.6923
$ws.Range("I5").Value2 = 865.375
$ws.Range("J5").Value2 = 608.8
$ws.Range("K5").Value2 = 2596.125
$ws.Range("L5").Value2 = 1826.4
$ws.Range("M5").Value2 = -2484.125
$ws.Range("N5").Value2 = -2050.4

# Row 86
$ws.Range("H86").Value2 = 711.25
$ws.Range("I86").Value2 = 679
$ws.Range("J86").Value2 = 765
$ws.Range("K86").Value2 = 2037
$ws.Range("L86").Value2 = 2295
$ws.Range("M86").Value2 = -851
$ws.Range("N86").Value2 = -4667

# Row 89
$ws.Range("H89").Value2 = 711.25
$ws.Range("I89").Value2 = 679
$ws.Range("J89").Value2 = 765
$ws.Range("K89").Value2 = 6111
$ws.Range("L89").Value2 = 6885
$ws.Range("M89").Value2 = -183
$ws.Range("N89").Value2 = -18741

# Row 97
$ws.Range("H97").Value2 = 371.85715
$ws.Range("I97").Value2 = 223.75
$ws.Range("J97").Value2 = 431.1
$ws.Range("K97").Value2 = 671.25
$ws.Range("L97").Value2 = 1293.3
$ws.Range("M97").Value2 = -175.25
$ws.Range("N97").Value2 = -2285.3

# Row 103
$ws.Range("H103").Value2 = 1610
$ws.Range("I103").Value2 = 810.4
$ws.Range("J103").Value2 = 1876.5333
$ws.Range("K103").Value2 = 2431.2
$ws.Range("L103").Value2 = 5629.5999
$ws.Range("M103").Value2 = -1552.2
$ws.Range("N103").Value2 = -7387.5999

# Row 114
$ws.Range("H114").Value2 = 892.1429000000001
$ws.Range("I114").Value2 = 804
$ws.Range("J114").Value2 = 1009.6667
$ws.Range("K114").Value2 = 2412
$ws.Range("L114").Value2 = 3029.0001
$ws.Range("M114").Value2 = 842
$ws.Range("N114").Value2 = -9537.000100000001

# Row 135
$ws.Range("H135").Value2 = 766.6923
$ws.Range("I135").Value2 = 865.375
$ws.Range("J135").Value2 = 608.8
$ws.Range("K135").Value2 = 7788.375
$ws.Range("L135").Value2 = 5479.2
$ws.Range("M135").Value2 = -5253.375
$ws.Range("N135").Value2 = -10549.2

# Row 136
$ws.Range("H136").Value2 = 2611.111
$ws.Range("I136").Value2 = 2611.111
$ws.Range("K136").Value2 = 7833.333
$ws.Range("M136").Value2 = -2733.333

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value2 = 4665.1665
$ws.Range("I43").Value2 = 4665.1665
$ws.Range("K43").Value2 = 4665.1665
$ws.Range("M43").Value2 = -4514.1665

# Row 80
$ws.Range("H80").Value2 = 2384.2703
$ws.Range("I80").Value2 = 1885.75
$ws.Range("K80").Value2 = 1885.75
$ws.Range("M80").Value2 = -887.75

# Row 83
$ws.Range("H83").Value2 = 2384.2703
$ws.Range("I83").Value2 = 1885.75
$ws.Range("K83").Value2 = 9428.75
$ws.Range("M83").Value2 = -4436.75

# Row 122
$ws.Range("H122").Value2 = 1682.5714
$ws.Range("I122").Value2 = 1633.0303
$ws.Range("J122").Value2 = 2500
$ws.Range("K122").Value2 = 4899.090899999999
$ws.Range("L122").Value2 = 7500
$ws.Range("M122").Value2 = -2449.090899999999
$ws.Range("N122").Value2 = -12400

# Row 136
$ws.Range("H136").Value2 = 41808.5
$ws.Range("J136").Value2 = 41808.5
$ws.Range("L136").Value2 = 125425.5
$ws.Range("N136").Value2 = -130525.5

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value2 = 4067.762
$ws.Range("I100").Value2 = 3727.1667
$ws.Range("J100").Value2 = 4521.8887
$ws.Range("K100").Value2 = 3727.1667
$ws.Range("L100").Value2 = 4521.8887
$ws.Range("M100").Value2 = -3186.1667
$ws.Range("N100").Value2 = -5603.8887

# Row 132
$ws.Range("H132").Value2 = 43086.8
$ws.Range("I132").Value2 = 50388.56
$ws.Range("J132").Value2 = 6578
$ws.Range("K132").Value2 = 151165.68
$ws.Range("L132").Value2 = 19734
$ws.Range("M132").Value2 = -148635.68
$ws.Range("N132").Value2 = -24794

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value2 = 5450
$ws.Range("I96").Value2 = 5450
$ws.Range("K96").Value2 = 5450
$ws.Range("M96").Value2 = -4077

# Row 113
$ws.Range("H113").Value2 = 1166.6774
$ws.Range("I113").Value2 = 969.1852
$ws.Range("J113").Value2 = 2499.75
$ws.Range("K113").Value2 = 2907.5556
$ws.Range("L113").Value2 = 7499.25
$ws.Range("M113").Value2 = -737.5556000000001
$ws.Range("N113").Value2 = -11839.25

# Row 132
$ws.Range("H132").Value2 = 27062.162
$ws.Range("I132").Value2 = 31532.723
$ws.Range("J132").Value2 = 4070.7144
$ws.Range("K132").Value2 = 94598.16900000001
$ws.Range("L132").Value2 = 12212.1432
$ws.Range("M132").Value2 = -92068.16900000001
$ws.Range("N132").Value2 = -17272.1432

# Row 136
$ws.Range("H136").Value2 = 890
$ws.Range("I136").Value2 = 890
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 2670
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -120
$ws.Range("N136").ClearContents()
